# Generate Report for Handoff
# Updates the localization-status workbook: re-sorts/refreshes the three
# file rows (ffff67e95264..., ffffffc9bfdff1..., 9d9de0ab...) on the
# Overview/zh-cn/de-de sheets and marks the 9d9de0ab row "Ready for handoff"
# with refreshed handoff timestamps (new handoff cycle kicked off).

$wb = $excel.ActiveWorkbook

function Set-SheetValues {
    param($ws, $values)

    foreach ($addr in $values.Keys) {
        $ws.Range($addr).Value = $values[$addr]
    }

    foreach ($hl in $ws.Hyperlinks) {
        $cellRef = $hl.Range.Address().Replace('$', '')
        if ($values.ContainsKey($cellRef)) {
            $hl.TextToDisplay = $values[$cellRef]
        }
    }
}

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverviewValues = @{
        "A1" = 'File Name'
        "B1" = 'zh-cn'
        "C1" = 'de-de'
        "A2" = 'ffff67e95264-01fd-41ee-a1f9-4aceb53d1d84.md'
        "B2" = 'Handed back: in sync with en-US'
        "C2" = 'Handed back: in sync with en-US'
        "A3" = 'ffffffc9bfdff1-f276-467c-91d6-7fd8d26b2024.md'
        "B3" = 'Handed back: in sync with en-US'
        "C3" = 'Handed back: in sync with en-US'
        "A4" = '9d9de0ab-c32a-4b20-833a-dd98a9be1e2d.md'
        "B4" = 'Ready for handoff'
        "C4" = 'Ready for handoff'
        "A5" = '.localization-config'
        "B5" = 'Not to be localized'
        "C5" = 'Not to be localized'
    }
Set-SheetValues $wsOverview $wsOverviewValues

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCnValues = @{
        "A1" = 'Source File Name'
        "B1" = 'Status'
        "C1" = 'Latest Handoff File'
        "D1" = 'Latest Handoff Datetime'
        "E1" = 'Latest Target File'
        "F1" = 'Latest Handback File'
        "G1" = 'Latest Handback DateTime'
        "H1" = 'Handoff Reason'
        "I1" = 'Dependency From'
        "A2" = 'ffff67e95264-01fd-41ee-a1f9-4aceb53d1d84.md'
        "B2" = 'Handed back: in sync with en-US'
        "C2" = '5e667d6f-dc90-457e-b295-fbe6abeb0028.3821691a11d655d326ebf1527d572680262db9ec.zh-cn.xlf'
        "D2" = '2016-01-28 05:51:49'
        "E2" = '5e667d6f-dc90-457e-b295-fbe6abeb0028.md'
        "F2" = '5e667d6f-dc90-457e-b295-fbe6abeb0028.3821691a11d655d326ebf1527d572680262db9ec.zh-cn.xlf'
        "G2" = '2016-01-28 05:52:30'
        "H2" = 'Include'
        "A3" = 'ffffffc9bfdff1-f276-467c-91d6-7fd8d26b2024.md'
        "B3" = 'Handed back: in sync with en-US'
        "C3" = '5e667d6f-dc90-457e-b295-fbe6abeb0028.3821691a11d655d326ebf1527d572680262db9ec.zh-cn.xlf'
        "D3" = '2016-01-28 05:51:49'
        "E3" = '5e667d6f-dc90-457e-b295-fbe6abeb0028.md'
        "F3" = '5e667d6f-dc90-457e-b295-fbe6abeb0028.3821691a11d655d326ebf1527d572680262db9ec.zh-cn.xlf'
        "G3" = '2016-01-28 05:52:30'
        "H3" = 'Include'
        "A4" = '9d9de0ab-c32a-4b20-833a-dd98a9be1e2d.md'
        "B4" = 'Ready for handoff'
        "C4" = '9d9de0ab-c32a-4b20-833a-dd98a9be1e2d.0e107ebc9410dd48340fd46ec6a53911ae5d41c3.zh-cn.xlf'
        "D4" = '2016-01-28 05:55:50'
        "E4" = '9d9de0ab-c32a-4b20-833a-dd98a9be1e2d.md'
        "F4" = '9d9de0ab-c32a-4b20-833a-dd98a9be1e2d.0e107ebc9410dd48340fd46ec6a53911ae5d41c3.zh-cn.xlf'
        "G4" = '2016-01-28 05:54:47'
        "H4" = 'Include'
        "A5" = '.localization-config'
        "B5" = 'Not to be localized'
        "D5" = '0001-01-01 00:00:00'
        "G5" = '0001-01-01 00:00:00'
        "H5" = 'Ignored'
    }
Set-SheetValues $wsZhCn $wsZhCnValues

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDeValues = @{
        "A1" = 'Source File Name'
        "B1" = 'Status'
        "C1" = 'Latest Handoff File'
        "D1" = 'Latest Handoff Datetime'
        "E1" = 'Latest Target File'
        "F1" = 'Latest Handback File'
        "G1" = 'Latest Handback DateTime'
        "H1" = 'Handoff Reason'
        "I1" = 'Dependency From'
        "A2" = 'ffff67e95264-01fd-41ee-a1f9-4aceb53d1d84.md'
        "B2" = 'Handed back: in sync with en-US'
        "C2" = '5e667d6f-dc90-457e-b295-fbe6abeb0028.3821691a11d655d326ebf1527d572680262db9ec.de-de.xlf'
        "D2" = '2016-01-28 05:51:59'
        "E2" = '5e667d6f-dc90-457e-b295-fbe6abeb0028.md'
        "F2" = '5e667d6f-dc90-457e-b295-fbe6abeb0028.3821691a11d655d326ebf1527d572680262db9ec.de-de.xlf'
        "G2" = '2016-01-28 05:52:47'
        "H2" = 'Include'
        "A3" = 'ffffffc9bfdff1-f276-467c-91d6-7fd8d26b2024.md'
        "B3" = 'Handed back: in sync with en-US'
        "C3" = '5e667d6f-dc90-457e-b295-fbe6abeb0028.3821691a11d655d326ebf1527d572680262db9ec.de-de.xlf'
        "D3" = '2016-01-28 05:51:59'
        "E3" = '5e667d6f-dc90-457e-b295-fbe6abeb0028.md'
        "F3" = '5e667d6f-dc90-457e-b295-fbe6abeb0028.3821691a11d655d326ebf1527d572680262db9ec.de-de.xlf'
        "G3" = '2016-01-28 05:52:47'
        "H3" = 'Include'
        "A4" = '9d9de0ab-c32a-4b20-833a-dd98a9be1e2d.md'
        "B4" = 'Ready for handoff'
        "C4" = '9d9de0ab-c32a-4b20-833a-dd98a9be1e2d.0e107ebc9410dd48340fd46ec6a53911ae5d41c3.de-de.xlf'
        "D4" = '2016-01-28 05:56:00'
        "E4" = '9d9de0ab-c32a-4b20-833a-dd98a9be1e2d.md'
        "F4" = '9d9de0ab-c32a-4b20-833a-dd98a9be1e2d.0e107ebc9410dd48340fd46ec6a53911ae5d41c3.de-de.xlf'
        "G4" = '2016-01-28 05:55:07'
        "H4" = 'Include'
        "A5" = '.localization-config'
        "B5" = 'Not to be localized'
        "D5" = '0001-01-01 00:00:00'
        "G5" = '0001-01-01 00:00:00'
        "H5" = 'Ignored'
    }
Set-SheetValues $wsDeDe $wsDeDeValues
